# Remove the "hospital outpatient clinic facility" row (BCIO:026015).
# This deletes row 28 entirely, shifting all rows below it up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(28).Delete()
